$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-7
# from 45170 (2023-09-01) to 45174 (2023-09-05)
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 3).Value = 45174
}
